$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The header cell "Month" should become lowercase "month".
$ws.Range("A1").Value = "month"

# Move the active selection back to A1 (the saved file had a stray
# selection at F4; Excel resets it to the default top-left cell).
[void]$ws.Range("A1").Select()
